$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2-9 down to 3-10)
$ws.Rows.Item(2).Insert()

# Fill in the new row with the new election data
$ws.Range("A2").Value = "Presidencia Municipal 15"
$ws.Range("B2").Value = "pm_15"
$ws.Range("C2").Value = "#669bbc"

# Update selection to match the new active cell/selection
$ws.Range("A2:C2").Select()
